$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Ativacao date: 01/01/2020 -> 01/01/2022 (row 8, cols B & C)
#    Assigning a date-shaped string straight to .Value lets Excel's
#    autoconvert turn it into a date serial, so we stage the text in
#    a scratch cell formatted as Text, then copy only the VALUE over
#    to B8/C8 so their original styles are left untouched.
# ------------------------------------------------------------------
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"
$scratch.Value = "01/01/2022"
$scratch.Copy()
$ws.Range("B8").PasteSpecial(-4163)
$ws.Range("C8").PasteSpecial(-4163)
$scratch.Clear()

# ------------------------------------------------------------------
# 2) Programa resumido (row 14): drop "Dupla Filtracao; "
# ------------------------------------------------------------------
$ws.Range("B14").Value = "Sistemas de Abastecimento de Água; Tecnologias de Tratamento de Água; Tratamento de Água em Ciclo Completo; Desinfecção; Filtração Direta Ascendente; Filtração Direta Descendente; Floto-Filtração; Filtração em Múltiplas Etapas; Tratamento dos Resíduos Gerados nas Estações de Tratamento de Água."
$ws.Range("C14").Value = "Sistemas de Abastecimento de Água; Tecnologias de Tratamento de Água; Tratamento de Água em Ciclo Completo; Desinfecção; Filtração Direta Ascendente; Filtração Direta Descendente; Floto-Filtração; Filtração em Múltiplas Etapas; Tratamento dos Resíduos Gerados nas Estações de Tratamento de Água."

# ------------------------------------------------------------------
# 3) Short syllabus (row 15): drop "Double Filtration; "
# ------------------------------------------------------------------
$ws.Range("B15").Value = "Systems of Water Supply; Water Treatment Technologies; Water Treatment in Complete Cycle; Disinfection, Direct Ascendant Filtration; Direct Descendant Filtration, Floto-filtration; Filtration in Multiple Levels, Treatment of the Generated Waste in the Water Treatment Stations."
$ws.Range("C15").Value = "Systems of Water Supply; Water Treatment Technologies; Water Treatment in Complete Cycle; Disinfection, Direct Ascendant Filtration; Direct Descendant Filtration, Floto-filtration; Filtration in Multiple Levels, Treatment of the Generated Waste in the Water Treatment Stations."

# ------------------------------------------------------------------
# 4) Programa (row 16): full rewrite of the detailed PT program
# ------------------------------------------------------------------
$ws.Range("B16").Value = "- Tipos de dimensionamento de redes de distribuição de água;- Características das águas de interesse para o tratamento: características físicas, químicas e bacteriológicas; - Padrão de Potabilidade; - Tecnologias de Tratamento de Água;- Unidades Constituintes de um Sistema de Abastecimento de Água;- Captação de Águas Subterrâneas e Captação de Águas Superficiais - Gradeamento, remoção de areia, casa de bombas; - Reservação; - Redes de distribuição: tipos de rede, peças e órgãos acessórios; - Sistema de Tratamento de Água de Ciclo Completo; - Coagulação-floculação e Mistura Rápida; - Decantação: decantação convencional e de alta taxa e sistema de remoção de lodo;- Mecanismos da filtração, materiais filtrantes e fundos de filtros, hidráulica da filtração, filtração com taxa constante e taxa declinante, - Desinfecção: principais desinfetantes, cloração e cloro-amoniação, pré e pós-cloração, - Tratamento de resíduos gerados nas ETAs e reuso de água recuperada"
$ws.Range("C16").Value = "- Tipos de dimensionamento de redes de distribuição de água;- Características das águas de interesse para o tratamento: características físicas, químicas e bacteriológicas; - Padrão de Potabilidade; - Tecnologias de Tratamento de Água;- Unidades Constituintes de um Sistema de Abastecimento de Água;- Captação de Águas Subterrâneas e Captação de Águas Superficiais - Gradeamento, remoção de areia, casa de bombas; - Reservação; - Redes de distribuição: tipos de rede, peças e órgãos acessórios; - Sistema de Tratamento de Água de Ciclo Completo; - Coagulação-floculação e Mistura Rápida; - Decantação: decantação convencional e de alta taxa e sistema de remoção de lodo;- Mecanismos da filtração, materiais filtrantes e fundos de filtros, hidráulica da filtração, filtração com taxa constante e taxa declinante, - Desinfecção: principais desinfetantes, cloração e cloro-amoniação, pré e pós-cloração, - Tratamento de resíduos gerados nas ETAs e reuso de água recuperada"

# ------------------------------------------------------------------
# 5) Syllabus (row 17): full rewrite of the detailed EN syllabus
# ------------------------------------------------------------------
$ws.Range("B17").Value = "- Types of dimensioning of water distribution networks;- Water characteristics of interest for treatment: physical, chemical and bacteriological characteristics;- Potability Standard;- Water Treatment Technologies;- Constituent Units of a Water Supply System;- Groundwater Catchment and Surface Water Catchment- Railing, sand removal, pump room;- Reservation;- Distribution networks: types of network, parts and accessories;- Full Cycle Water Treatment System;- Coagulation-flocculation and Rapid Mixing;- Decantation: conventional and high rate decantation and sludge removal system;- Filtration mechanisms, filter materials and filter bottoms, filtration hydraulics, filtration with constant rate and declining rate,- Disinfection: main disinfectants, chlorination and chlor-ammonia, pre and post-chlorination,- Treatment of waste generated at stations and reuse of recovered water"
$ws.Range("C17").Value = "- Types of dimensioning of water distribution networks;- Water characteristics of interest for treatment: physical, chemical and bacteriological characteristics;- Potability Standard;- Water Treatment Technologies;- Constituent Units of a Water Supply System;- Groundwater Catchment and Surface Water Catchment- Railing, sand removal, pump room;- Reservation;- Distribution networks: types of network, parts and accessories;- Full Cycle Water Treatment System;- Coagulation-flocculation and Rapid Mixing;- Decantation: conventional and high rate decantation and sludge removal system;- Filtration mechanisms, filter materials and filter bottoms, filtration hydraulics, filtration with constant rate and declining rate,- Disinfection: main disinfectants, chlorination and chlor-ammonia, pre and post-chlorination,- Treatment of waste generated at stations and reuse of recovered water"

# ------------------------------------------------------------------
# 6) Metodo (row 19): trim down to just the evaluation sentence
# ------------------------------------------------------------------
$ws.Range("B19").Value = "Avaliação baseada em provas, exercícios, trabalhos práticos e relatórios."
$ws.Range("C19").Value = "Avaliação baseada em provas, exercícios, trabalhos práticos e relatórios."
